# Apply the updated cryptocurrency price/volume snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.132.98"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "1.578.49"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "'212.09"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("E6").Value = "  +6.70%  "
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").Value = "'25.93"
$ws.Range("E8").Value = "  +9.67%  "
$ws.Range("E9").Value = "  +2.44%  "
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("D11").Value = "'0.0900"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "1.803.32"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").Value = "1.569.54"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "29.155.76"
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'3.70"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").Value = "'62.33"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").Value = "'238.81"
$ws.Range("E18").Value = "  +5.23%  "
$ws.Range("D19").Value = "'7.44"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").Value = "0.0₃0692"
$ws.Range("E20").Value = "  +2.75%  "
$ws.Range("D21").Value = "'0.997"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "'3.99"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").Value = "'9.20"
$ws.Range("E23").Value = "  +4.21%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  +5.11%  "
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("E26").Value = "  +4.37%  "
$ws.Range("D27").Value = "'15.14"
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").Value = "'0.0466"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("D33").Value = "1.424.75"
$ws.Range("E33").Value = "  +3.11%  "
$ws.Range("D34").Value = "'3.06"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("D37").Value = "'2.76"
$ws.Range("E37").Value = "  +6.90%  "
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").Value = "'0.527"
$ws.Range("E40").Value = "  +3.63%  "
$ws.Range("D41").Value = "'1.96"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("D42").Value = "'53.73"
$ws.Range("E42").Value = "  +27.22%  "
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "'0.789"
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("D45").Value = "'0.0470"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").Value = "'64.45"
$ws.Range("E46").Value = "  +4.30%  "
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "1.716.21"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("E49").Value = "  -6.52%  "
$ws.Range("D50").Value = "'85.77"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "'0.0512"
$ws.Range("E51").Value = "  +0.71%  "
